# Auto-generated edit script: refresh BRVM "Recommandations" and "Top_YTD" tables
$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsTop  = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations": columns A..G, rows 2..44 ---
# row, Titre, JoursHausse, JoursBaisse, VariationTotale, DerniereVariation, Recommandation, Strategie
$recoData = @(
    @(2, "NEI-CEDA CI", 0, 4, 3730, 900, "🟡 Observer", "➖ Neutre"),
    @(3, "BRVM - SERVICES PUBLICS", 0, 8, 3397.16, 111.85, "🟡 Observer", "➖ Neutre"),
    @(4, "BRVM - AUTRES SECTEURS", 0, 4, 2355.65, 581.47, "🟡 Observer", "➖ Neutre"),
    @(5, "BRVM - DISTRIBUTION", 0, 4, 2029.85, 503.78, "🟡 Observer", "➖ Neutre"),
    @(6, "BRVM - TRANSPORT", 0, 4, 1438.04, 359.82, "🟡 Observer", "➖ Neutre"),
    @(7, "BRVM - AGRICULTURE", 0, 4, 1353.5, 339.1, "🟡 Observer", "➖ Neutre"),
    @(8, "BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 4, 681.66, 168.17, "🟡 Observer", "➖ Neutre"),
    @(9, "BRVM - FINANCES", 0, 4, 589.7, 147.41, "🟡 Observer", "➖ Neutre"),
    @(10, "BRVM - SERVICES FINANCIERS", 0, 4, 579.5599999999999, 144.88, "🟡 Observer", "➖ Neutre"),
    @(11, "BRVM-PRESTIGE", 0, 4, 575.0599999999999, 143, "🟡 Observer", "➖ Neutre"),
    @(12, "BRVM - INDUSTRIE                 (**)", 0, 2, 530.45, 266.96, "🟡 Observer", "➖ Neutre"),
    @(13, "BRVM - INDUSTRIELS", 0, 4, 498.56, 125.48, "🟡 Observer", "➖ Neutre"),
    @(14, "BRVM - ENERGIE", 0, 4, 458.61, 114.87, "🟡 Observer", "➖ Neutre"),
    @(15, "BRVM - CONSOMMATION DE BASE          (**)", 0, 2, 444.46, 222.9, "🟡 Observer", "➖ Neutre"),
    @(16, "BRVM-PRINCIPAL                    (**)", 0, 2, 439.09, 219.74, "🟡 Observer", "➖ Neutre"),
    @(17, "BRVM - TELECOMMUNICATIONS", 0, 4, 383.99, 95.43000000000001, "🟡 Observer", "➖ Neutre"),
    @(18, "BRVM - INDUSTRIE", 0, 1, 264.92, 264.92, "🟡 Observer", "➖ Neutre"),
    @(19, "BRVM - CONSOMMATION DE BASE         (**)", 0, 1, 220.71, 220.71, "🟡 Observer", "➖ Neutre"),
    @(20, "BRVM-PRINCIPAL                   (**)", 0, 1, 220.21, 220.21, "🟡 Observer", "➖ Neutre"),
    @(21, "ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)", 4, 0, 29.83, 7.41, "🟢 Achat", "✅ Renforcer"),
    @(22, "SERVAIR ABIDJAN CI (ABJC)", 3, 0, 22.2, 7.37, "🟢 Achat", "✅ Renforcer"),
    @(23, "SUCRIVOIRE (SCRC)", 2, 1, 7.65, 7.5, "🟡 Observer", "👀 À surveiller"),
    @(24, "ECOBANK TRANS. INCORP. TG (ETIT)", 1, 0, 4.76, 4.76, "🟡 Observer", "➖ Neutre"),
    @(25, "SICABLE CI (CABC)", 1, 0, 4.55, 4.55, "🟡 Observer", "➖ Neutre"),
    @(26, "SOLIBRA CI (SLBC)", 2, 1, 2.98, 7.49, "🟡 Observer", "👀 À surveiller"),
    @(27, "BERNABE CI (BNBC)", 1, 1, 2.14, 7.35, "🟡 Observer", "👀 À surveiller"),
    @(28, "AFRICA GLOBAL LOGISTICS CI (SDSC)", 1, 0, 2.07, 2.07, "🟡 Observer", "➖ Neutre"),
    @(29, "BANK OF AFRICA ML (BOAM)", 1, 0, 1.96, 1.96, "🟡 Observer", "➖ Neutre"),
    @(30, "ONATEL BF (ONTBF)", 1, 0, 1.81, 1.81, "🟡 Observer", "➖ Neutre"),
    @(31, "CORIS BANK INTERNATIONAL (CBIBF)", 1, 0, 1.43, 1.43, "🟡 Observer", "➖ Neutre"),
    @(32, "BANK OF AFRICA BN (BOAB)", 1, 1, 0.08, 2.76, "🟡 Observer", "👀 À surveiller"),
    @(33, "TOTAL", 0, 4, 0, 0, "🟡 Observer", "➖ Neutre"),
    @(34, "ORANGE COTE D'IVOIRE (ORAC)", 0, 1, -1.3, -1.3, "🟡 Observer", "➖ Neutre"),
    @(35, "NESTLE CI (NTLC)", 0, 1, -1.32, -1.32, "🟡 Observer", "➖ Neutre"),
    @(36, "SONATEL SN (SNTS)", 0, 1, -1.96, -1.96, "🟡 Observer", "➖ Neutre"),
    @(37, "SICOR CI (SICC)", 0, 1, -2.5, -2.5, "🟡 Observer", "➖ Neutre"),
    @(38, "FILTISAC CI (FTSC)", 0, 1, -2.83, -2.83, "🟡 Observer", "➖ Neutre"),
    @(39, "ORAGROUP TOGO (ORGT)", 0, 1, -3.1, -3.1, "🟡 Observer", "➖ Neutre"),
    @(40, "ECOBANK COTE D''IVOIRE (ECOC)", 0, 2, -5.63, -3.73, "🟡 Observer", "➖ Neutre"),
    @(41, "SAFCA CI (SAFC)", 0, 1, -6.09, -6.09, "🟡 Observer", "➖ Neutre"),
    @(42, "UNILEVER CI (UNLC)", 0, 1, -6.97, -6.97, "🟡 Observer", "➖ Neutre"),
    @(43, "SETAO CI (STAC)", 0, 3, -9.77, -1.87, "🔴 Vente", "⚠️ Risque de décrochage"),
    @(44, "CFAO MOTORS CI (CFAC)", 0, 2, -10.38, -5.94, "🟡 Observer", "➖ Neutre")
)

foreach ($row in $recoData) {
    $r = $row[0]
    $wsReco.Cells.Item($r, 1).Value = $row[1]
    $wsReco.Cells.Item($r, 2).Value = $row[2]
    $wsReco.Cells.Item($r, 3).Value = $row[3]
    $wsReco.Cells.Item($r, 4).Value = $row[4]
    $wsReco.Cells.Item($r, 5).Value = $row[5]
    $wsReco.Cells.Item($r, 6).Value = $row[6]
    $wsReco.Cells.Item($r, 7).Value = $row[7]
}

# --- Sheet "Top_YTD": columns A..B, rows 2..11 ---
# row, Titre, ProgressionYTD
$topData = @(
    @(2, "BRVM - SERVICES PUBLICS", 9919240.6),
    @(3, "NEI-CEDA CI", 1135606),
    @(4, "BRVM - AUTRES SECTEURS", 225124.74),
    @(5, "BRVM - DISTRIBUTION", 136054.67),
    @(6, "BRVM - TRANSPORT", 44480.8),
    @(7, "BRVM - AGRICULTURE", 36830.12),
    @(8, "BRVM - CONSOMMATION DISCRETIONNAIRE", 5246.36),
    @(9, "BRVM - FINANCES", 3647.77),
    @(10, "BRVM - SERVICES FINANCIERS", 3496.52),
    @(11, "BRVM-PRESTIGE", 3430.85)
)

foreach ($row in $topData) {
    $r = $row[0]
    $wsTop.Cells.Item($r, 1).Value = $row[1]
    $wsTop.Cells.Item($r, 2).Value = $row[2]
}

